$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "phone Case Type" column (D) and the "tracking Number" column
# (originally S, which becomes R once D is removed). This shifts every
# column to the right of D left by one, and removes the trailing tracking
# Number column.
$ws.Columns("D").Delete()
$ws.Columns("R").Delete()

# --- Update remaining cell values to match the new data set ---

# Row 1 (headers) are unchanged by the above column deletes.

# Row 2
$ws.Range("A2").Value = 43750
$ws.Range("D2").Value = "dr1001"

# Row 3
$ws.Range("A3").Value = 43751
$ws.Range("D3").Value = "dr1002"

# Row 4
$ws.Range("A4").Value = 43752
$ws.Range("D4").Value = "dr1003"

# Row 5
$ws.Range("A5").Value = 43753
$ws.Range("D5").Value = "dr1004"

# --- Restore the view state ---
$ws.Range("E11").Select()
